# Add "Hint" and "Popup" columns (H, I) to the TwoWays process sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("H1").Value = "Hint"
$ws.Range("I1").Value = "Popup"

$headerRange = $ws.Range("H1:I1")

# Match the look of the existing header cells (F1/G1): white text on a
# blue solid fill, General number format, no border, no special
# alignment/wrap - but with a distinct (swapped) blue fill color.
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.Font.Color = 16777215
$headerRange.Interior.PatternColor = 10506806
$headerRange.Interior.Color = 10773812

# Update the active selection to match the edited workbook state.
$ws.Range("E15").Select()
